$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.236.89"
$ws.Range("E2").Value = "  -1.49%  "
$ws.Range("D3").Value = "2.514.27"
$ws.Range("E3").Value = "  -0.26%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'571.85"
$ws.Range("E5").Value = "  -0.59%  "
$ws.Range("D6").Value = "'165.32"
$ws.Range("E6").Value = "  -2.61%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "'0.514"
$ws.Range("E8").Value = "  +0.82%  "
$ws.Range("D9").Value = "2.514.06"
$ws.Range("E9").Value = "  -0.19%  "
$ws.Range("E10").Value = "  +0.38%  "
$ws.Range("D11").Value = "'0.168"
$ws.Range("E11").Value = "  -0.27%  "
$ws.Range("E12").Value = "  +3.57%  "
$ws.Range("D13").Value = "'4.91"
$ws.Range("E13").Value = "  +1.80%  "
$ws.Range("D14").Value = "2.975.96"
$ws.Range("E14").Value = "  -0.18%  "
$ws.Range("D15").Value = "69.111.66"
$ws.Range("E15").Value = "  -1.30%  "
$ws.Range("E16").Value = "  -2.41%  "
$ws.Range("D17").Value = "'24.72"
$ws.Range("E17").Value = "  -1.09%  "
$ws.Range("D18").Value = "2.520.95"
$ws.Range("E18").Value = "  +0.14%  "
$ws.Range("D19").Value = "'11.28"
$ws.Range("E19").Value = "  -1.88%  "
$ws.Range("D20").Value = "'7.65"
$ws.Range("E20").Value = "  +0.94%  "
$ws.Range("D21").Value = "'347.81"
$ws.Range("E21").Value = "  -1.95%  "
$ws.Range("E22").Value = "  -0.94%  "
$ws.Range("E23").Value = "  +0.23%  "
$ws.Range("E24").Value = "  -0.13%  "
$ws.Range("D25").Value = "'70.28"
$ws.Range("E25").Value = "  +1.84%  "
$ws.Range("D26").Value = "'3.92"
$ws.Range("E26").Value = "  -4.20%  "
$ws.Range("D27").Value = "'8.85"
$ws.Range("E27").Value = "  -3.79%  "
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("D29").Value = "'0.998"
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("D30").Value = "0.0₃0888"
$ws.Range("E30").Value = "  -2.44%  "
$ws.Range("D31").Value = "'7.78"
$ws.Range("E31").Value = "  -0.86%  "
$ws.Range("D32").Value = "'459.33"
$ws.Range("E32").Value = "  -4.63%  "
$ws.Range("E33").Value = "  -4.61%  "
$ws.Range("E34").Value = "  -1.71%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.14%  "
$ws.Range("E36").Value = "  -0.06%  "
$ws.Range("D37").Value = "'157.27"
$ws.Range("E37").Value = "  +0.69%  "
$ws.Range("D38").Value = "'19.01"
$ws.Range("E38").Value = "  +0.66%  "
$ws.Range("D39").Value = "'18.46"
$ws.Range("E39").Value = "  -0.56%  "
$ws.Range("E40").Value = "  -0.04%  "
$ws.Range("E41").Value = "  -0.79%  "
$ws.Range("D42").Value = "'4.69"
$ws.Range("E42").Value = "  -0.88%  "
$ws.Range("E43").Value = "  -3.06%  "
$ws.Range("D44").Value = "'38.08"
$ws.Range("E44").Value = "  -0.56%  "
$ws.Range("E45").Value = "  -8.31%  "
$ws.Range("E46").Value = "  -7.35%  "
$ws.Range("D47").Value = "'141.68"
$ws.Range("E47").Value = "  -0.42%  "
$ws.Range("E48").Value = "  -1.01%  "
$ws.Range("D49").Value = "'3.46"
$ws.Range("E49").Value = "  -1.72%  "
$ws.Range("D50").Value = "'0.0729"
$ws.Range("E50").Value = "  -0.31%  "
$ws.Range("E51").Value = "  -3.33%  "
